$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: status cells for zh-cn / de-de columns
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# zh-cn sheet: status, handback datetime, error detail cleared
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-08-26 00:48:33"
$wsZhCn.Range("P2").Value = "'"
$wsZhCn.Range("P2").Style = "Normal"
$wsZhCn.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZhCn.Columns.Item(16).ColumnWidth = 13.7470528738839

# de-de sheet: status, handback datetime, error detail cleared
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-26 00:48:39"
$wsDeDe.Range("P2").Value = "'"
$wsDeDe.Range("P2").Style = "Normal"
$wsDeDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDeDe.Columns.Item(16).ColumnWidth = 13.7470528738839
